# ---------------------------------------------------------------------------
# Adds a new "2022-Q3" quarterly sheet (right after "总计") with its fund
# holdings table, and updates the "总计" (summary) sheet with the new
# 2022-Q3 row, shifting the rest of the history down by one row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Update the "总计" summary sheet: insert the 2022-Q3 row at the top of
#    the data and push the existing history rows down by one.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

$zongji = @(
    ,@(2, 0, "2022-Q3", 31, 1.41)
    ,@(3, 1, "2022-Q2", 43, 3.6)
    ,@(4, 2, "2022-Q1", 30, 2.9)
    ,@(5, 3, "2021-Q4", 50, 7.56)
    ,@(6, 4, "2021-Q3", 62, 13.59)
    ,@(7, 5, "2021-Q2", 41, 9.2)
    ,@(8, 6, "2021-Q1", 64, 14.82)
    ,@(9, 7, "2020-Q4", 34, 8.94)
)

foreach ($row in $zongji) {
    $r = $row[0]
    $summary.Cells.Item($r, 1).Value = $row[1]
    $summary.Cells.Item($r, 2).Value = $row[2]
    $summary.Cells.Item($r, 3).Value = $row[3]
    $summary.Cells.Item($r, 4).Value = $row[4]
}

# Row 9 is brand new (the sheet used to stop at row 8) - give A9 the same
# look (bold, centered, bordered) as the rest of the index column by
# copying the format from A8.
$summary.Cells.Item(8, 1).Copy()
$summary.Cells.Item(9, 1).PasteSpecial(-4122)
$summary.Cells.Item(9, 1).Value = 7
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 2) Insert the new "2022-Q3" worksheet right after "总计" and populate it
#    with the fund holdings table for that quarter.
# ---------------------------------------------------------------------------
$q3sheet = $wb.Worksheets.Add($null, $summary)
$q3sheet.Name = "2022-Q3"

# Header row - reuse the bold/centered/bordered header style (style index 2)
# from the summary sheet's header row.
$summary.Range("B1").Copy()
$q3sheet.Range("B1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q3sheet.Range("B1").Value = "基金代码"
$q3sheet.Range("C1").Value = "基金名称"
$q3sheet.Range("D1").Value = "基金规模"
$q3sheet.Range("E1").Value = "股票总仓位"
$q3sheet.Range("F1").Value = "仓位占比"
$q3sheet.Range("G1").Value = "持有市值(亿元)"
$q3sheet.Range("H1").Value = "仓位排名"

# Column A (row index) also uses the bold/centered/bordered style (style
# index 2), same as the summary sheet's index column.
$summary.Range("A2").Copy()
$q3sheet.Range("A2:A32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$q3 = @(
    ,@(0, "003713", "英大睿盛灵活配置混合A", "2.83", "93.65", "6.32", "0.1789", 7)
    ,@(1, "001305", "九泰天富改革新动力混合A", "1.86", "94.71", "9.18", "0.1707", 2)
    ,@(2, "167508", "安信价值发现两年定期开放混合（LOF）", "3.26", "89.75", "5.13", "0.1672", 6)
    ,@(3, "003714", "英大睿盛灵活配置混合C", "2.19", "93.65", "6.32", "0.1384", 7)
    ,@(4, "008704", "广发高股息优享混合A", "2.33", "92.85", "4.96", "0.1156", 9)
    ,@(5, "003345", "安信新成长灵活配置混合A", "4.98", "32.19", "1.71", "0.0852", 6)
    ,@(6, "001844", "九泰久益灵活配置混合C", "0.98", "93.32", "8.50", "0.0833", 3)
    ,@(7, "003029", "安信新优选灵活配置混合C", "3.61", "33.25", "2.20", "0.0794", 4)
    ,@(8, "004138", "上银鑫达灵活配置混合A", "1.86", "75.90", "3.50", "0.0651", 8)
    ,@(9, "001782", "九泰久益灵活配置混合A", "0.53", "93.32", "8.50", "0.0450", 3)
    ,@(10, "001607", "英大策略优选混合A", "0.57", "91.98", "6.80", "0.0388", 4)
    ,@(11, "008705", "广发高股息优享混合C", "0.69", "92.85", "4.96", "0.0342", 9)
    ,@(12, "009766", "安信平稳双利3个月持有期混合A", "0.92", "30.99", "2.94", "0.0270", 4)
    ,@(13, "004249", "安信中国制造混合", "0.52", "89.55", "5.06", "0.0263", 6)
    ,@(14, "004393", "安信企业价值优选混合", "0.54", "87.43", "4.49", "0.0242", 7)
    ,@(15, "007393", "上银未来生活灵活配置混合A", "0.74", "81.08", "3.02", "0.0223", 7)
    ,@(16, "003447", "英大睿鑫灵活配置混合C", "0.21", "92.71", "7.88", "0.0165", 3)
    ,@(17, "012522", "英大稳固增强核心一年持有混合C", "1.24", "27.71", "1.32", "0.0164", 9)
    ,@(18, "001399", "安信鑫安得利灵活配置混合A", "0.70", "30.15", "2.30", "0.0161", 5)
    ,@(19, "009912", "九泰天富改革新动力混合C", "0.17", "94.71", "9.18", "0.0156", 2)
    ,@(20, "012521", "英大稳固增强核心一年持有混合A", "0.75", "27.71", "1.32", "0.0099", 9)
    ,@(21, "003346", "安信新成长灵活配置混合C", "0.46", "32.19", "1.71", "0.0079", 6)
    ,@(22, "001400", "安信鑫安得利灵活配置混合C", "0.25", "30.15", "2.30", "0.0058", 5)
    ,@(23, "003446", "英大睿鑫灵活配置混合A", "0.07", "92.71", "7.88", "0.0055", 3)
    ,@(24, "009767", "安信平稳双利3个月持有期混合C", "0.17", "30.99", "2.94", "0.0050", 4)
    ,@(25, "014113", "上银未来生活灵活配置混合C", "0.16", "81.08", "3.02", "0.0048", 7)
    ,@(26, "750005", "安信平稳增长混合A", "0.07", "58.57", "3.72", "0.0026", 7)
    ,@(27, "001608", "英大策略优选混合C", "0.02", "91.98", "6.80", "0.0014", 4)
    ,@(28, "003028", "安信新优选灵活配置混合A", "0.06", "33.25", "2.20", "0.0013", 4)
    ,@(29, "002035", "安信平稳增长混合C", "0.03", "58.57", "3.72", "0.0011", 7)
    ,@(30, "015753", "上银鑫达灵活配置混合C", "0.03", "75.90", "3.50", "0.0010", 8)
)

# Columns B-G hold numeric-looking text (fund code, name, size, position...)
# that must stay TEXT (e.g. "003713", "8.50"), so force text formatting on
# the whole block before writing the values.
$q3sheet.Range("B2:G32").NumberFormat = "@"

$r = 2
foreach ($row in $q3) {
    $q3sheet.Cells.Item($r, 1).Value = $row[0]
    $q3sheet.Cells.Item($r, 2).Value = $row[1]
    $q3sheet.Cells.Item($r, 3).Value = $row[2]
    $q3sheet.Cells.Item($r, 4).Value = $row[3]
    $q3sheet.Cells.Item($r, 5).Value = $row[4]
    $q3sheet.Cells.Item($r, 6).Value = $row[5]
    $q3sheet.Cells.Item($r, 7).Value = $row[6]
    $q3sheet.Cells.Item($r, 8).Value = $row[7]
    $r++
}

$summary.Activate()

Write-Output "2022-Q3 sheet added and 总计 updated"
